# Add a fixed course ("fix" / "a132" / SaturdayT1 everyWeek) to every course
# block already present in the schedule output sheet.
#
# Each existing course block looks like:
#   row N   (bold/shaded): ID, Penalty
#   row N+1              : Course name, Teacher name, Times   <- data row
#   row N+2              : (blank separator)
#
# We keep the existing data row in place, insert a fresh row right below it
# to hold the course that used to live there, and overwrite the original
# data row with the new fixed course's details. Working from the bottom of
# the sheet upward keeps the not-yet-processed row numbers stable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRows = @(24, 21, 18, 15, 12, 9, 6, 3)

foreach ($r in $dataRows) {
    $courseName = $ws.Cells.Item($r, 3).Value()
    $teacherName = $ws.Cells.Item($r, 4).Value()
    $times = $ws.Cells.Item($r, 5).Value()

    $ws.Rows.Item($r + 1).Insert()

    $ws.Cells.Item($r + 1, 3).Value = $courseName
    $ws.Cells.Item($r + 1, 4).Value = $teacherName
    $ws.Cells.Item($r + 1, 5).Value = $times

    $ws.Cells.Item($r, 3).Value = "fix"
    $ws.Cells.Item($r, 4).Value = "a132"
    $ws.Cells.Item($r, 5).Value = "SaturdayT1  everyWeek  "
}
